$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Spent" amount for Rent
$ws.Range("C2").Value = 10000

# Insert a new row at 6 (above "Fun") for the new "Coffee" category
$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = "Coffee"
$ws.Range("B6").Value = 500
$ws.Range("C6").Value = 480

# Insert a new row at 8 (below the shifted "Fun" row, above "Tithe") for "Gas"
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "Gas"
$ws.Range("B8").Value = 1000
$ws.Range("C8").Value = 1600

# The "Total" row (now at row 11, with a blank row 10 above it) needs its
# formulas to cover the expanded data range B2:B9 / C2:C9
$ws.Range("B11").Formula = "=SUM(B2:B9)"
$ws.Range("C11").Formula = "=SUM(C2:C9)"

# Restore the active selection shown in the saved workbook
$ws.Range("E7").Select()
